$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 716 (shifts existing rows 716..757 down to 717..758)
$ws.Rows("716").Insert()

# Populate the new row with the inserted daily record: 2026/01/24, 土, 19, 154
# Force the date-like text to be stored as plain text (not auto-converted to a
# date serial number) by temporarily switching the cell to Text format, then
# clear the explicit formatting afterwards so the cell keeps the default
# (unstyled) look, matching the rest of the data rows in the sheet.
$ws.Range("A716").NumberFormat = "@"
$ws.Range("A716").Value = "2026/01/24"
$ws.Range("A716").ClearFormats()

$ws.Range("B716").Value = "土"
$ws.Range("C716").Value = 19
$ws.Range("D716").Value = 154
